$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Process.time (Seconds) column (B2:B11) to 2 for all rows
$ws.Range("B2:B11").Value = 2

# Update the active cell selection to B9
$ws.Range("B9").Select()
